$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price (D) and Volume(1h) (E) columns store plain text values that
# often look numeric (e.g. "23.297.93", "0.00001187", "297.91"). Force
# Text format on exactly the cells being rewritten so Excel keeps the
# literal string instead of parsing/rounding it as a number.
$ws.Range("D2:D3").NumberFormat = "@"
$ws.Range("D5:D11").NumberFormat = "@"
$ws.Range("D13:D44").NumberFormat = "@"
$ws.Range("D46:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '23.297.93'
$ws.Range("E2").Value = '  -1.91%  '
$ws.Range("D3").Value = '1.627.32'
$ws.Range("E3").Value = '  -1.99%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = '1.002'
$ws.Range("E5").Value = '  +0.14%  '
$ws.Range("D6").Value = '297.91'
$ws.Range("E6").Value = '  -1.66%  '
$ws.Range("D7").Value = '0.3763'
$ws.Range("E7").Value = '  -1.57%  '
$ws.Range("D8").Value = '49.82'
$ws.Range("E8").Value = '  -2.70%  '
$ws.Range("D9").Value = '0.3470'
$ws.Range("E9").Value = '  -4.12%  '
$ws.Range("D10").Value = '0.08022'
$ws.Range("E10").Value = '  -2.12%  '
$ws.Range("D11").Value = '1.196'
$ws.Range("E11").Value = '  -2.91%  '
$ws.Range("E12").Value = '  +0.20%  '
$ws.Range("D13").Value = '21.84'
$ws.Range("E13").Value = '  -3.29%  '
$ws.Range("D14").Value = '6.276'
$ws.Range("E14").Value = '  -3.06%  '
$ws.Range("D15").Value = '7.204'
$ws.Range("E15").Value = '  -3.05%  '
$ws.Range("D16").Value = '0.00001187'
$ws.Range("E16").Value = '  -3.29%  '
$ws.Range("D17").Value = '1.627.61'
$ws.Range("E17").Value = '  -1.98%  '
$ws.Range("D18").Value = '94.51'
$ws.Range("E18").Value = '  -3.43%  '
$ws.Range("D19").Value = '0.06939'
$ws.Range("E19").Value = '  -1.18%  '
$ws.Range("D20").Value = '6.593'
$ws.Range("E20").Value = '  -3.39%  '
$ws.Range("D21").Value = '17.22'
$ws.Range("E21").Value = '  -2.45%  '
$ws.Range("D22").Value = '1.002'
$ws.Range("E22").Value = '  +0.20%  '
$ws.Range("D23").Value = '12.32'
$ws.Range("E23").Value = '  -4.00%  '
$ws.Range("D24").Value = '23.312.90'
$ws.Range("E24").Value = '  -1.85%  '
$ws.Range("D25").Value = '2.410'
$ws.Range("E25").Value = '  -3.84%  '
$ws.Range("D26").Value = '2.947'
$ws.Range("E26").Value = '  -1.83%  '
$ws.Range("D27").Value = '20.80'
$ws.Range("E27").Value = '  -1.97%  '
$ws.Range("D28").Value = '150.17'
$ws.Range("E28").Value = '  -1.71%  '
$ws.Range("D29").Value = '5.157'
$ws.Range("E29").Value = '  -1.41%  '
$ws.Range("D30").Value = '130.33'
$ws.Range("E30").Value = '  -2.91%  '
$ws.Range("D31").Value = '1.812.51'
$ws.Range("E31").Value = '  -1.61%  '
$ws.Range("D32").Value = '6.700'
$ws.Range("E32").Value = '  -6.12%  '
$ws.Range("D33").Value = '2.119'
$ws.Range("E33").Value = '  -5.38%  '
$ws.Range("D34").Value = '11.16'
$ws.Range("E34").Value = '  -7.41%  '
$ws.Range("D35").Value = '0.9740'
$ws.Range("E35").Value = '  -7.80%  '
$ws.Range("D36").Value = '0.02648'
$ws.Range("E36").Value = '  -6.00%  '
$ws.Range("D37").Value = '0.08736'
$ws.Range("E37").Value = '  -0.90%  '
$ws.Range("D38").Value = '0.2405'
$ws.Range("E38").Value = '  -4.47%  '
$ws.Range("D39").Value = '5.804'
$ws.Range("E39").Value = '  -4.72%  '
$ws.Range("D40").Value = '0.06713'
$ws.Range("E40").Value = '  -4.36%  '
$ws.Range("D41").Value = '12.62'
$ws.Range("E41").Value = '  -3.11%  '
$ws.Range("D42").Value = '0.6774'
$ws.Range("E42").Value = '  -3.35%  '
$ws.Range("D43").Value = '1.287'
$ws.Range("E43").Value = '  -3.62%  '
$ws.Range("D44").Value = '15.27'
$ws.Range("E44").Value = '  -4.68%  '
$ws.Range("E45").Value = '  +0.14%  '
$ws.Range("D46").Value = '0.6285'
$ws.Range("E46").Value = '  -3.65%  '
$ws.Range("D47").Value = '2.222'
$ws.Range("E47").Value = '  -3.71%  '
$ws.Range("D48").Value = '3.883'
$ws.Range("E48").Value = '  -2.10%  '
$ws.Range("D49").Value = '0.07606'
$ws.Range("E49").Value = '  -3.95%  '
$ws.Range("D50").Value = '125.90'
$ws.Range("E50").Value = '  -2.03%  '
$ws.Range("D51").Value = '1.213'
$ws.Range("E51").Value = '  +1.50%  '
